# Fruta / hortaliza, semanal
# Update weekly price records: dates and volume/price figures are refreshed
# for several rows (rows 4 and 7 are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value  = 44189   # D2 Fecha
$ws.Cells.Item(2, 13).Value = 40      # M2 Volumen
$ws.Cells.Item(2, 14).Value = 15000   # N2 Precio minimo
$ws.Cells.Item(2, 15).Value = 15000   # O2 Precio maximo
$ws.Cells.Item(2, 16).Value = 15000   # P2 Precio promedio ponderado
$ws.Cells.Item(2, 19).Value = 3000    # S2 Precio $/Kg

# Row 3
$ws.Cells.Item(3, 4).Value  = 44188   # D3 Fecha
$ws.Cells.Item(3, 13).Value = 30      # M3 Volumen

# Row 5
$ws.Cells.Item(5, 4).Value  = 44931   # D5 Fecha
$ws.Cells.Item(5, 13).Value = 50      # M5 Volumen
$ws.Cells.Item(5, 14).Value = 18000   # N5 Precio minimo
$ws.Cells.Item(5, 15).Value = 18000   # O5 Precio maximo
$ws.Cells.Item(5, 16).Value = 18000   # P5 Precio promedio ponderado
$ws.Cells.Item(5, 19).Value = 3600    # S5 Precio $/Kg

# Row 6
$ws.Cells.Item(6, 4).Value  = 44902   # D6 Fecha
$ws.Cells.Item(6, 13).Value = 35      # M6 Volumen
$ws.Cells.Item(6, 14).Value = 12000   # N6 Precio minimo
$ws.Cells.Item(6, 15).Value = 12000   # O6 Precio maximo
$ws.Cells.Item(6, 16).Value = 12000   # P6 Precio promedio ponderado
$ws.Cells.Item(6, 19).Value = 2400    # S6 Precio $/Kg

# Row 8
$ws.Cells.Item(8, 4).Value  = 44196   # D8 Fecha
$ws.Cells.Item(8, 13).Value = 56      # M8 Volumen

# Row 9
$ws.Cells.Item(9, 4).Value  = 44907   # D9 Fecha
$ws.Cells.Item(9, 13).Value = 45      # M9 Volumen
$ws.Cells.Item(9, 14).Value = 25000   # N9 Precio minimo
$ws.Cells.Item(9, 15).Value = 25000   # O9 Precio maximo
$ws.Cells.Item(9, 16).Value = 25000   # P9 Precio promedio ponderado
$ws.Cells.Item(9, 19).Value = 5000    # S9 Precio $/Kg

# Row 10
$ws.Cells.Item(10, 4).Value  = 44193  # D10 Fecha
$ws.Cells.Item(10, 13).Value = 40     # M10 Volumen
$ws.Cells.Item(10, 14).Value = 15000  # N10 Precio minimo
$ws.Cells.Item(10, 15).Value = 15000  # O10 Precio maximo
$ws.Cells.Item(10, 16).Value = 15000  # P10 Precio promedio ponderado
$ws.Cells.Item(10, 19).Value = 3000   # S10 Precio $/Kg

# Row 11
$ws.Cells.Item(11, 4).Value  = 44179  # D11 Fecha
$ws.Cells.Item(11, 13).Value = 45     # M11 Volumen
$ws.Cells.Item(11, 14).Value = 20000  # N11 Precio minimo
$ws.Cells.Item(11, 15).Value = 20000  # O11 Precio maximo
$ws.Cells.Item(11, 16).Value = 20000  # P11 Precio promedio ponderado
$ws.Cells.Item(11, 19).Value = 4000   # S11 Precio $/Kg

# Row 12
$ws.Cells.Item(12, 4).Value  = 44914  # D12 Fecha
$ws.Cells.Item(12, 13).Value = 56     # M12 Volumen
$ws.Cells.Item(12, 14).Value = 23000  # N12 Precio minimo
$ws.Cells.Item(12, 15).Value = 23000  # O12 Precio maximo
$ws.Cells.Item(12, 16).Value = 23000  # P12 Precio promedio ponderado
$ws.Cells.Item(12, 19).Value = 4600   # S12 Precio $/Kg

$wb.Save()
